$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inflation")
Write-Host $ws.Name
